# Re-process the sheet with the newly curated dimensions.
# Most columns that used to be "iaest-dimension:*" metric columns are now
# proper "iaest-measure:*" columns; only the reference-period (ano) and the
# two area-name columns (municipio-nombre / provincia-nombre) remain real
# dimensions. Column N (municipio-nombre) is promoted from a plain
# "iaest-measure:municipio-nombre" measure to the curated
# "sdmx-dimension:refArea" dimension, matching provincia-nombre (column P).
# Because these columns are no longer "dim" mapped against an external
# workbook, their row-5 mapping-file references are removed - except for
# "ano" (column B), which stays a curated dimension with its own mapping file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: concept URI ---------------------------------------------------
$ws.Range("E2").Value = "iaest-measure:salario-perceptores"
$ws.Range("G2").Value = "iaest-measure:ccaa"
$ws.Range("H2").Value = "iaest-measure:salario-medio-anual"
$ws.Range("K2").Value = "iaest-measure:pension-percepciones-por-persona"
$ws.Range("N2").Value = "sdmx-dimension:refArea"
$ws.Range("Q2").Value = "iaest-measure:pension-perceptores"
$ws.Range("R2").Value = "iaest-measure:desempleo-medio-por-percepcion"
$ws.Range("T2").Value = "iaest-measure:pension-media-por-persona"
$ws.Range("U2").Value = "iaest-measure:salario-percepciones-por-persona"
$ws.Range("V2").Value = "iaest-measure:salario-medio-por-persona"

# --- Row 3: medida / dim ---------------------------------------------------
$ws.Range("E3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("K3").Value = "medida"
$ws.Range("N3").Value = "dim"
$ws.Range("Q3").Value = "medida"
$ws.Range("R3").Value = "medida"
$ws.Range("T3").Value = "medida"
$ws.Range("U3").Value = "medida"
$ws.Range("V3").Value = "medida"

# --- Row 4: datatype / URI template ---------------------------------------
$ws.Range("E4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("K4").Value = "xsd:int"
$ws.Range("N4").Value = "URI-Municipio"
$ws.Range("Q4").Value = "xsd:int"
$ws.Range("R4").Value = "xsd:int"
$ws.Range("T4").Value = "xsd:int"
$ws.Range("U4").Value = "xsd:int"
$ws.Range("V4").Value = "xsd:int"

# --- Row 5: mapping workbook - no longer needed for former dimensions -----
$ws.Range("E5").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("H5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("Q5").ClearContents()
$ws.Range("R5").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("U5").ClearContents()
$ws.Range("V5").ClearContents()
